$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: Normal-flow step 9 sentence in the "Login" use-case table.
#   "9. System navigate to the user information page."
#     -> "9. System navigate to the "Course list" page."
# The old sentence carried the (hidden) "_GoBack" bookmark between
# "the" and " user information page." - that bookmark is relocated
# below (Change 2), so here we simply replace the tail of the
# sentence, which also removes the bookmark from this spot.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("the user information page.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "the "
    $rng.Collapse(0)
    $rng.InsertAfter([char]0x201C + "Course list" + [char]0x201D + " page.")
}

# ------------------------------------------------------------------
# Change 2: In the Alternative-Flow A3 pop-up text, split "e-mail"
# into "E" / "-mail" and drop the relocated "_GoBack" bookmark right
# between them:
#   ...alert pop-up text that "e-mail or Password does not exist".
#   -> ...alert pop-up text that "E[_GoBack]-mail or Password does not exist".
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("pop-up text that " + [char]0x201C + "e-mail", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    [void]$rng2.MoveStart(1, 18)
    [void]$rng2.MoveEnd(1, -5)
    $rng2.Text = "E"
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("pop-up text that " + [char]0x201C + "E-mail", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    [void]$rng3.MoveStart(1, 19)
    $rng3.Collapse(1)
    [void]$d.Bookmarks.Add("_GoBack", $rng3)
}
